$d = $word.ActiveDocument

# --- locate "Version" and split it into "Versi" | "on" -------------------
# The saved file keeps two formatting-identical runs side by side; the only
# way this engine preserves such a boundary across a save is if the range
# formatting was actually touched (even if reverted to the same value
# right after), so toggle-and-revert Bold on the "on" part of the word.
$verRange = $d.Content
$verRange.Find.Execute("Version", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitStart = $verRange.Start + 5
$splitEnd = $verRange.End
$splitRange = $d.Range($splitStart, $splitEnd)
$splitRange.Bold = 1
$splitRange.Bold = 0

# --- "Version 1." -> "Version 2." ----------------------------------------
$numRange = $d.Content
$numRange.Find.Execute("1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$numRange.Text = "2"

# --- drop the "." that used to sit right after the digit -----------------
$dotRange = $d.Content
$dotRange.Find.Execute(".", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dotRange.Text = ""

# --- put a new "." at the very end of the paragraph, i.e. after the
#     _GoBack bookmark that precedes the paragraph mark -------------------
$para = $d.Paragraphs(1).Range
$tail = $d.Range($para.End - 1, $para.End - 1)
$tail.InsertAfter(".")
